# Apply crypto price/volume updates per commit "Updated cryptos list on Thu Jan 18 18:09:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'41.736.65"
$ws.Range("E2").Value = "  -1.66%  "

# Row 3
$ws.Range("D3").Value = "'2.473.40"
$ws.Range("E3").Value = "  -2.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "'310.64"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("D6").Value = "'94.93"
$ws.Range("E6").Value = "  -4.40%  "

# Row 7
$ws.Range("D7").Value = "'0.553"
$ws.Range("E7").Value = "  -2.65%  "

# Row 8
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$ws.Range("D9").Value = "'0.509"
$ws.Range("E9").Value = "  -3.79%  "

# Row 10
$ws.Range("D10").Value = "'33.91"
$ws.Range("E10").Value = "  -5.30%  "

# Row 11
$ws.Range("E11").Value = "  -2.53%  "

# Row 12
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("D13").Value = "'6.97"
$ws.Range("E13").Value = "  -4.96%  "

# Row 14
$ws.Range("D14").Value = "'2.861.74"
$ws.Range("E14").Value = "  -2.03%  "

# Row 15
$ws.Range("D15").Value = "'2.528.22"
$ws.Range("E15").Value = "  -3.14%  "

# Row 16
$ws.Range("D16").Value = "'14.65"
$ws.Range("E16").Value = "  -6.93%  "

# Row 17
$ws.Range("D17").Value = "'0.789"
$ws.Range("E17").Value = "  -4.12%  "

# Row 18
$ws.Range("D18").Value = "'41.741.24"
$ws.Range("E18").Value = "  -1.62%  "

# Row 19
$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = "  -6.29%  "

# Row 20
$ws.Range("E20").Value = "  -3.47%  "

# Row 21
$ws.Range("D21").Value = "'11.64"
$ws.Range("E21").Value = "  -4.85%  "

# Row 22
$ws.Range("D22").Value = "'69.54"
$ws.Range("E22").Value = "  +0.45%  "

# Row 23
$ws.Range("D23").Value = "'236.07"
$ws.Range("E23").Value = "  -3.24%  "

# Row 24
$ws.Range("D24").Value = "'2.78"
$ws.Range("E24").Value = "  -3.50%  "

# Row 25
$ws.Range("E25").Value = "  -5.25%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$ws.Range("D27").Value = "'24.70"
$ws.Range("E27").Value = "  -4.79%  "

# Row 28
$ws.Range("E28").Value = "  -4.84%  "

# Row 29
$ws.Range("E29").Value = "  -3.72%  "

# Row 30
$ws.Range("D30").Value = "'36.32"
$ws.Range("E30").Value = "  -7.50%  "

# Row 31
$ws.Range("D31").Value = "'154.39"
$ws.Range("E31").Value = "  -2.17%  "

# Row 32
$ws.Range("D32").Value = "'5.62"
$ws.Range("E32").Value = "  -1.85%  "

# Row 33
$ws.Range("D33").Value = "'2.63"
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("B34").Value = "ApeXProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "  -7.79%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0755"
$ws.Range("E35").Value = "  -4.98%  "

# Row 36
$ws.Range("E36").Value = "  -4.19%  "

# Row 37
$ws.Range("D37").Value = "'17.23"
$ws.Range("E37").Value = "  -5.49%  "

# Row 38
$ws.Range("E38").Value = "  -6.85%  "

# Row 39
$ws.Range("D39").Value = "'0.106"
$ws.Range("E39").Value = "  -5.04%  "

# Row 40
$ws.Range("E40").Value = "  -3.59%  "

# Row 41
$ws.Range("E41").Value = "  -6.67%  "

# Row 42
$ws.Range("D42").Value = "'21.22"
$ws.Range("E42").Value = "  -3.03%  "

# Row 43
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").Value = "'1.993.71"
$ws.Range("E44").Value = "  +1.60%  "

# Row 45
$ws.Range("E45").Value = "  -3.80%  "

# Row 46
$ws.Range("E46").Value = "  -7.16%  "

# Row 47
$ws.Range("D47").Value = "'8.64"
$ws.Range("E47").Value = "  -2.88%  "

# Row 48
$ws.Range("D48").Value = "'2.712.39"
$ws.Range("E48").Value = "  -2.26%  "

# Row 49
$ws.Range("D49").Value = "'76.41"
$ws.Range("E49").Value = "  -5.66%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'97.79"
$ws.Range("E50").Value = "  -3.63%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.181"
$ws.Range("E51").Value = "  -5.59%  "

